$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.115.99"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "2.565.64"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'586.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.10%  "
$ws.Range("D6").Value = "'148.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("E9").Value = "  +2.69%  "
$ws.Range("E10").Value = "  +1.16%  "
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "'0.356"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").Value = "'27.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").Value = "3.028.82"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Value = "62.976.71"
$ws.Range("E16").Value = "  +2.07%  "
$ws.Range("D17").Value = "2.562.70"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").Value = "'343.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("E20").Value = "  +2.93%  "
$ws.Range("D21").Value = "'6.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("E23").Value = "  -3.65%  "
$ws.Range("D24").Value = "'66.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("D25").Value = "2.679.62"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").Value = "'1.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("E28").Value = "  +11.46%  "
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'8.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "'1.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.31%  "
$ws.Range("D33").Value = "0.0₃0827"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("D34").Value = "'461.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.98%  "
$ws.Range("D35").Value = "'176.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").Value = "'1.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.62%  "
$ws.Range("D37").Value = "'0.405"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("D39").Value = "'4.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.98%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").Value = "'1.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "'151.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").Value = "'0.0550"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.02%  "
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("D51").Value = "'11.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.66%  "
